# Apply updated dSF (column F) values for the specified rows.
# Commit message: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -1
    13 = -1
    17 = -3
    18 = -1
    20 = 1
    26 = 2
    27 = -2
    29 = -2
    32 = 0
    33 = -3
    35 = 1
    42 = -2
    46 = 3
    49 = 1
    51 = 3
    55 = -2
    59 = 4
    66 = -1
    68 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
